# Clava Parser Restructuring - update node-completion status markers
# across the decl/expr/stmt/type/attr/other sheets, refresh the Summary
# rollup formula, add a new "ClavaNode" entry to the "other" sheet, and
# update sheet selections / active tab to match the author's final state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Formula = "=(B3 + (B4/2))/B2"

# ---------------------------------------------------------------
# decl sheet: clear two old markers, add two new ones
# ---------------------------------------------------------------
$decl = $wb.Worksheets.Item("decl")
$decl.Range("B2").ClearContents()
$decl.Range("B3").ClearContents()
$decl.Range("B9").Value = "-"
$decl.Range("B11").Value = "o"

# ---------------------------------------------------------------
# expr sheet: add new markers
# ---------------------------------------------------------------
$expr = $wb.Worksheets.Item("expr")
$expr.Range("B6").Value = "-"
$expr.Range("B7").Value = "-"
$expr.Range("B12").Value = "o"
$expr.Range("B34").Value = "-"
$expr.Range("B36").Value = "-"
$expr.Range("B38").Value = "-"
$expr.Range("B40").Value = "-"
$expr.Range("B43").Value = "-"
$expr.Range("B45").Value = "-"
$expr.Range("B46").Value = "o"

# ---------------------------------------------------------------
# stmt sheet: add new markers
# ---------------------------------------------------------------
$stmt = $wb.Worksheets.Item("stmt")
$stmt.Range("B14").Value = "-"
$stmt.Range("B15").Value = "-"
$stmt.Range("B23").Value = "-"

# ---------------------------------------------------------------
# type sheet: add new markers
# ---------------------------------------------------------------
$type = $wb.Worksheets.Item("type")
$type.Range("B6").Value = "-"
$type.Range("B11").Value = "-"
$type.Range("B16").Value = "-"
$type.Range("B19").Value = "o"
$type.Range("B21").Value = "-"
$type.Range("B25").Value = "-"
$type.Range("B33").Value = "-"

# ---------------------------------------------------------------
# attr sheet: add new markers
# ---------------------------------------------------------------
$attr = $wb.Worksheets.Item("attr")
$attr.Range("B2").Value = "o"
$attr.Range("B3").Value = "o"
$attr.Range("B4").Value = "o"
$attr.Range("B5").Value = "-"
$attr.Range("B6").Value = "o"

# ---------------------------------------------------------------
# other sheet: insert a new "ClavaNode" row at the top of the data
# (shift column-A node names down by one row, keep the summary
# formulas anchored on row 2), then mark the existing "App" row.
# ---------------------------------------------------------------
$other = $wb.Worksheets.Item("other")
for ($r = 30; $r -ge 2; $r--) {
    $prev = $other.Cells.Item($r, 1).Value()
    $other.Cells.Item($r + 1, 1).Value = $prev
}
$other.Range("A2").Value = "ClavaNode"
$other.Range("B2").Value = "-"
$other.Range("B13").Value = "o"

# ---------------------------------------------------------------
# Restore per-sheet selections and active tab to match final state.
# The "type" sheet is activated last so it ends up as the selected tab.
# ---------------------------------------------------------------
$summary.Activate()
$summary.Range("B10").Select()

$decl.Activate()
$decl.Range("B12").Select()

$expr.Activate()
$expr.Range("B44").Select()

$stmt.Activate()
$stmt.Range("B16").Select()

$attr.Activate()
$attr.Range("B5").Select()

$other.Activate()
$other.Range("B14").Select()

$type.Activate()
$type.Range("R20").Select()
